{"js": "// Replace each two-digit-division answer cell's text with its new value.\n// Each `old` string is unique in the document, so body.search(...) returns\n// exactly one hit per pair. We gather ALL search hits first (before any\n// writes land) so a newly-written value can never be re-matched by a later\n// search in this same batch (e.g. '63\u00f73=21, 0' -> '28\u00f79=3, 1' and separately\n// '28\u00f79=3, 1' -> '50\u00f79=5, 5').\nconst body = context.document.body;\nconst pairs = [\n  [\"98\u00f75=19, 3\", \"14\u00f77=2, 0\"],\n  [\"62\u00f73=20, 2\", \"91\u00f77=13, 0\"],\n  [\"29\u00f72=14, 1\", \"37\u00f74=9, 1\"],\n  [\"16\u00f75=3, 1\", \"97\u00f75=19, 2\"],\n  [\"84\u00f73=28, 0\", \"27\u00f78=3, 3\"],\n  [\"63\u00f73=21, 0\", \"28\u00f79=3, 1\"],\n  [\"64\u00f79=7, 1\", \"87\u00f78=10, 7\"],\n  [\"32\u00f73=10, 2\", \"84\u00f78=10, 4\"],\n  [\"32\u00f76=5, 2\", \"60\u00f72=30, 0\"],\n  [\"29\u00f74=7, 1\", \"77\u00f72=38, 1\"],\n  [\"36\u00f76=6, 0\", \"84\u00f74=21, 0\"],\n  [\"28\u00f74=7, 0\", \"53\u00f78=6, 5\"],\n  [\"88\u00f73=29, 1\", \"88\u00f76=14, 4\"],\n  [\"80\u00f76=13, 2\", \"33\u00f76=5, 3\"],\n  [\"67\u00f79=7, 4\", \"13\u00f78=1, 5\"],\n  [\"35\u00f78=4, 3\", \"82\u00f79=9, 1\"],\n  [\"59\u00f77=8, 3\", \"87\u00f76=14, 3\"],\n  [\"82\u00f74=20, 2\", \"69\u00f78=8, 5\"],\n  [\"99\u00f75=19, 4\", \"36\u00f75=7, 1\"],\n  [\"28\u00f79=3, 1\", \"50\u00f79=5, 5\"],\n  [\"58\u00f74=14, 2\", \"36\u00f75=7, 1\"],\n  [\"71\u00f76=11, 5\", \"83\u00f74=20, 3\"],\n  [\"39\u00f78=4, 7\", \"20\u00f75=4, 0\"],\n  [\"23\u00f79=2, 5\", \"39\u00f72=19, 1\"],\n  [\"31\u00f75=6, 1\", \"48\u00f77=6, 6\"],\n];\n\nconst searchResults = pairs.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((r) => r.load(\"items\"));\nawait context.sync();\n\nsearchResults.forEach((r, i) => {\n  const [oldText, newText] = pairs[i];\n  if (r.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n  r.items[0].insertText(newText, \"Replace\");\n});\nawait context.sync();\n", "ps1": "# Update each two-digit-division answer cell in the single 20x5 table with\n# its new value. Cells are targeted by (row, column) rather than by\n# text search, since several old/new values collide across cells (e.g. one\n# cell's new text equals another cell's old text) which would make a naive\n# sequential Find/Replace overwrite the wrong cell.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$updates = @(\n  @{ Row = 1; Col = 1; OldText = \"98\u00f75=19, 3\"; NewText = \"14\u00f77=2, 0\" }\n  @{ Row = 1; Col = 2; OldText = \"62\u00f73=20, 2\"; NewText = \"91\u00f77=13, 0\" }\n  @{ Row = 1; Col = 3; OldText = \"29\u00f72=14, 1\"; NewText = \"37\u00f74=9, 1\" }\n  @{ Row = 1; Col = 4; OldText = \"16\u00f75=3, 1\"; NewText = \"97\u00f75=19, 2\" }\n  @{ Row = 1; Col = 5; OldText = \"84\u00f73=28, 0\"; NewText = \"27\u00f78=3, 3\" }\n  @{ Row = 5; Col = 1; OldText = \"63\u00f73=21, 0\"; NewText = \"28\u00f79=3, 1\" }\n  @{ Row = 5; Col = 2; OldText = \"64\u00f79=7, 1\"; NewText = \"87\u00f78=10, 7\" }\n  @{ Row = 5; Col = 3; OldText = \"32\u00f73=10, 2\"; NewText = \"84\u00f78=10, 4\" }\n  @{ Row = 5; Col = 4; OldText = \"32\u00f76=5, 2\"; NewText = \"60\u00f72=30, 0\" }\n  @{ Row = 5; Col = 5; OldText = \"29\u00f74=7, 1\"; NewText = \"77\u00f72=38, 1\" }\n  @{ Row = 9; Col = 1; OldText = \"36\u00f76=6, 0\"; NewText = \"84\u00f74=21, 0\" }\n  @{ Row = 9; Col = 2; OldText = \"28\u00f74=7, 0\"; NewText = \"53\u00f78=6, 5\" }\n  @{ Row = 9; Col = 3; OldText = \"88\u00f73=29, 1\"; NewText = \"88\u00f76=14, 4\" }\n  @{ Row = 9; Col = 4; OldText = \"80\u00f76=13, 2\"; NewText = \"33\u00f76=5, 3\" }\n  @{ Row = 9; Col = 5; OldText = \"67\u00f79=7, 4\"; NewText = \"13\u00f78=1, 5\" }\n  @{ Row = 13; Col = 1; OldText = \"35\u00f78=4, 3\"; NewText = \"82\u00f79=9, 1\" }\n  @{ Row = 13; Col = 2; OldText = \"59\u00f77=8, 3\"; NewText = \"87\u00f76=14, 3\" }\n  @{ Row = 13; Col = 3; OldText = \"82\u00f74=20, 2\"; NewText = \"69\u00f78=8, 5\" }\n  @{ Row = 13; Col = 4; OldText = \"99\u00f75=19, 4\"; NewText = \"36\u00f75=7, 1\" }\n  @{ Row = 13; Col = 5; OldText = \"28\u00f79=3, 1\"; NewText = \"50\u00f79=5, 5\" }\n  @{ Row = 17; Col = 1; OldText = \"58\u00f74=14, 2\"; NewText = \"36\u00f75=7, 1\" }\n  @{ Row = 17; Col = 2; OldText = \"71\u00f76=11, 5\"; NewText = \"83\u00f74=20, 3\" }\n  @{ Row = 17; Col = 3; OldText = \"39\u00f78=4, 7\"; NewText = \"20\u00f75=4, 0\" }\n  @{ Row = 17; Col = 4; OldText = \"23\u00f79=2, 5\"; NewText = \"39\u00f72=19, 1\" }\n  @{ Row = 17; Col = 5; OldText = \"31\u00f75=6, 1\"; NewText = \"48\u00f77=6, 6\" }\n)\n\nforeach ($u in $updates) {\n  $cell = $tbl.Cell($u.Row, $u.Col)\n  $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n  if ($current -ne $u.OldText) {\n    throw \"Cell ($($u.Row),$($u.Col)) expected '$($u.OldText)' but found '$current'\"\n  }\n  $cell.Range.Text = $u.NewText\n}\n"}
